# Update date and multiplication problems in the worksheet.

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-02-26 Monday"; new = "2024-02-27 Tuesday"},
    @{old = "392×6="; new = "813×3="},
    @{old = "176×4="; new = "762×4="},
    @{old = "851×6="; new = "534×4="},
    @{old = "619×8="; new = "180×8="},
    @{old = "469×7="; new = "245×2="},
    @{old = "209×7="; new = "798×5="},
    @{old = "885×7="; new = "660×2="},
    @{old = "547×7="; new = "905×4="},
    @{old = "917×2="; new = "362×7="},
    @{old = "264×8="; new = "647×8="},
    @{old = "680×3="; new = "356×2="},
    @{old = "361×8="; new = "659×8="},
    @{old = "809×4="; new = "740×2="},
    @{old = "494×3="; new = "149×7="},
    @{old = "259×8="; new = "972×6="},
    @{old = "504×9="; new = "144×4="},
    @{old = "591×6="; new = "227×5="},
    @{old = "352×6="; new = "127×7="},
    @{old = "185×5="; new = "373×6="},
    @{old = "875×2="; new = "229×4="},
    @{old = "396×8="; new = "509×8="},
    @{old = "897×6="; new = "257×2="},
    @{old = "550×5="; new = "487×6="},
    @{old = "501×3="; new = "700×2="},
    @{old = "864×4="; new = "144×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
